$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# Sheet 1: LP1912
# -----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 02:59:08"
$ws1.Range("A3").Value = "Total filas: 5"

# Row 6: 215_ALUAR
$ws1.Range("A6").Value = "02:59:08"
$ws1.Range("B6").Value = "02:59"
$ws1.Range("D6").Value = 0

# Row 7: 14_ABASTO
$ws1.Range("A7").Value = "02:59:08"
$ws1.Range("D7").Value = 49

# Row 8: 81_EL PELIGRO
$ws1.Range("A8").Value = "02:59:08"
$ws1.Range("D8").Value = 62

# New row 9: 215A_EL PATO
$ws1.Range("A9").Value = "02:59:08"
$ws1.Range("B9").Value = "04:46"
$ws1.Range("C9").Value = "215A_EL PATO"
$ws1.Range("D9").Value = 107
$ws1.Range("E9").Value = "LP1912"

# New row 10: 11_ETCHEVERRY
$ws1.Range("A10").Value = "02:59:08"
$ws1.Range("B10").Value = "04:53"
$ws1.Range("C10").Value = "11_ETCHEVERRY"
$ws1.Range("D10").Value = 114
$ws1.Range("E10").Value = "LP1912"

# -----------------------------------------------------------------------
# Sheet 2: LP1912-215
# -----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 02:59:08"
$ws2.Range("A3").Value = "Total filas: 2"

# Row 6: 215_ALUAR
$ws2.Range("A6").Value = "02:59:08"
$ws2.Range("B6").Value = "02:59"
$ws2.Range("D6").Value = 0

# New row 7: 215A_EL PATO
$ws2.Range("A7").Value = "02:59:08"
$ws2.Range("B7").Value = "04:46"
$ws2.Range("C7").Value = "215A_EL PATO"
$ws2.Range("D7").Value = 107
$ws2.Range("E7").Value = "LP1912"

# -----------------------------------------------------------------------
# Sheet 3: 6203-6173
# -----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 02:59:08"
